$wb = $excel.ActiveWorkbook

# Correct the "Web Data 49" placeholder text to "Web Data 1" in the header
# (cell C1) of every sheet - they all show the current Web Data / order
# period label.
foreach ($sheetName in @("OrderPeriod", "AddProductCategory1", "AddCustomer", "CreateOrder")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C1").Value = "Web Data 1"
}

# Added Validation for the order Period in the Product Dependency Feature:
# move the active selection from "CreateOrder" (cell I1) back to the
# "OrderPeriod" sheet (cell C1), which becomes the active/selected tab.
$createOrder = $wb.Worksheets.Item("CreateOrder")
$createOrder.Range("C1").Select() | Out-Null

$orderPeriod = $wb.Worksheets.Item("OrderPeriod")
$orderPeriod.Activate() | Out-Null
$orderPeriod.Range("C1").Select() | Out-Null
